$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.555.03"
$ws.Range("E2").Value = "  +1.99%  "

# Row 3
$ws.Range("D3").Value = "1.810.82"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.49%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.79"
$ws.Range("E5").Value = "  -2.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4509"
$ws.Range("E7").Value = "  +5.87%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3775"
$ws.Range("E8").Value = "  +7.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.86"
$ws.Range("E9").Value = "  -1.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.153"
$ws.Range("E10").Value = "  +0.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07556"
$ws.Range("E11").Value = "  +1.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.74"
$ws.Range("E12").Value = "  -1.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9987"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.328"
$ws.Range("E14").Value = "  +0.85%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.635"
$ws.Range("E15").Value = "  +4.18%  "

# Row 16
$ws.Range("D16").Value = "1.804.39"
$ws.Range("E16").Value = "  -0.80%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").Value = "  +0.57%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06774"
$ws.Range("E18").Value = "  +1.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.07"
$ws.Range("E19").Value = "  -1.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9987"
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.71"
$ws.Range("E21").Value = "  +2.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.337"
$ws.Range("E22").Value = "  -1.23%  "

# Row 23
$ws.Range("D23").Value = "28.533.42"
$ws.Range("E23").Value = "  +1.70%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.87"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.406"
$ws.Range("E25").Value = "  -0.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.70"
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.374"
$ws.Range("E27").Value = "  -3.96%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.66"
$ws.Range("E28").Value = "  -2.36%  "

# Row 29
$ws.Range("D29").Value = "2.006.55"
$ws.Range("E29").Value = "  -0.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.92"
$ws.Range("E30").Value = "  +0.89%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.272"
$ws.Range("E31").Value = "  -2.56%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.996"
$ws.Range("E32").Value = "  -1.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.859"
$ws.Range("E33").Value = "  -2.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09349"
$ws.Range("E34").Value = "  +2.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2273"
$ws.Range("E35").Value = "  +4.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.28"
$ws.Range("E36").Value = "  -0.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06374"
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02326"
$ws.Range("E38").Value = "  -1.95%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6622"
$ws.Range("E39").Value = "  -1.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.185"
$ws.Range("E40").Value = "  -1.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.216"
$ws.Range("E41").Value = "  -0.29%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.139"
$ws.Range("E42").Value = "  +0.61%  "

# Row 43
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.448"
$ws.Range("E43").Value = "  -3.54%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9983"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("E45").Value = "  -3.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6106"
$ws.Range("E46").Value = "  -0.73%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.813"
$ws.Range("E47").Value = "  -1.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.84"
$ws.Range("E48").Value = "  +0.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.044"
$ws.Range("E49").Value = "  -0.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07113"
$ws.Range("E50").Value = "  -0.28%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.160"
$ws.Range("E51").Value = "  -1.91%  "
